# The source export was re-sorted: the underlying sighting "records" for
# rows 12-13 (swap) and rows 15-18 (4-cycle rotation) moved to different
# row numbers while row 14 stayed put. Re-apply each record's field values
# to its new row, and drop the optional time/method cells that don't
# belong to the record now occupying that row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12 <- (was row 13's record) ---
$ws.Range("A12").Value  = 130983072
$ws.Range("Q12").Value  = 570809
$ws.Range("R12").Value  = 6736404
$ws.Range("S12").Value  = 10
$ws.Range("Z12").Value  = "08:44"
$ws.Range("AB12").Value = "08:44"
$ws.Range("AF12").Value = ""
$ws.Range("AW12").Value = "Bo karlstens"
$ws.Range("AX12").Value = "Bo karlstens"

# --- Row 13 <- (was row 12's record) ---
$ws.Range("A13").Value  = 130979103
$ws.Range("Q13").Value  = 570739
$ws.Range("R13").Value  = 6736418
$ws.Range("S13").Value  = 1
$ws.Range("AW13").Value = "Erik Danielsson"
$ws.Range("AX13").Value = "Erik Danielsson"
$ws.Range("Z13").ClearContents()
$ws.Range("AB13").ClearContents()
$ws.Range("AF13").ClearContents()

# --- Row 15 <- (was row 18's record) ---
$ws.Range("A15").Value  = 130979082
$ws.Range("B15").Value  = 57884
$ws.Range("E15").Value  = 100109
$ws.Range("F15").Value  = "Tretåig hackspett"
$ws.Range("G15").Value  = "Picoides tridactylus"
$ws.Range("H15").Value  = "(Linnaeus, 1758)"
$ws.Range("Q15").Value  = 570952
$ws.Range("R15").Value  = 6736563
$ws.Range("S15").Value  = 1
$ws.Range("AC15").Value = "Äldre ringhack"
$ws.Range("AW15").Value = "Erik Danielsson"
$ws.Range("AX15").Value = "Erik Danielsson"
$ws.Range("Z15").ClearContents()
$ws.Range("AB15").ClearContents()
$ws.Range("AF15").ClearContents()

# --- Row 16 <- (was row 15's record) ---
$ws.Range("A16").Value  = 130983071
$ws.Range("P16").Value  = "Flytjärnsmyran, Dlr"
$ws.Range("Q16").Value  = 570817
$ws.Range("R16").Value  = 6736417
$ws.Range("Z16").Value  = "08:53"
$ws.Range("AB16").Value = "08:53"
$ws.Range("AF16").Value = ""
$ws.Range("AW16").Value = "Bo karlstens"
$ws.Range("AX16").Value = "Bo karlstens"

# --- Row 17 <- (was row 16's record) ---
$ws.Range("A17").Value  = 130983619
$ws.Range("P17").Value  = "Flytjärnsmyren, Dlr"
$ws.Range("Q17").Value  = 570825
$ws.Range("R17").Value  = 6736389
$ws.Range("Z17").Value  = "08:54"
$ws.Range("AB17").Value = "08:54"
$ws.Range("AW17").Value = "Göran Ehn"
$ws.Range("AX17").Value = "Göran Ehn"
$ws.Range("AF17").ClearContents()

# --- Row 18 <- (was row 17's record) ---
$ws.Range("A18").Value  = 130983074
$ws.Range("B18").Value  = 79244
$ws.Range("E18").Value  = 6425
$ws.Range("F18").Value  = "Garnlav"
$ws.Range("G18").Value  = "Alectoria sarmentosa"
$ws.Range("H18").Value  = "(Ach.) Ach."
$ws.Range("Q18").Value  = 570764
$ws.Range("R18").Value  = 6736425
$ws.Range("S18").Value  = 10
$ws.Range("Z18").Value  = "08:23"
$ws.Range("AB18").Value = "08:23"
$ws.Range("AF18").Value = ""
$ws.Range("AW18").Value = "Bo karlstens"
$ws.Range("AX18").Value = "Bo karlstens"
$ws.Range("AC18").ClearContents()
